# fix (bug #2): Error indicador 1.4
#
# The biomass generation indicator (criterion C1.4) in the alternative_info
# sheet was computed with a bug; this corrects the raw indicator values
# (alternative_info!L), their normalization (alternatives_norm!D) and the
# downstream AHP evaluation scores (result!B). Also bumps the run date
# recorded on the info sheet.

$wb = $excel.ActiveWorkbook

$wsInfo    = $wb.Worksheets.Item("info")
$wsAltInfo = $wb.Worksheets.Item("alternative_info")
$wsAltNorm = $wb.Worksheets.Item("alternatives_norm")
$wsResult  = $wb.Worksheets.Item("result")

# info sheet - run date
# Leading apostrophe keeps this a literal text value (matches the original
# inlineStr cell) instead of having COM auto-parse "04/02/23" into a date.

$wsInfo.Range("C2").Formula = "'04/02/23"

# alternative_info sheet - raw C1.4 indicator values (col L) recomputed
$wsAltInfo.Range("L5").Value = 0.958603951173407
$wsAltInfo.Range("L7").Value = 0.8510959963766491
$wsAltInfo.Range("L8").Value = 1.131928037118369
$wsAltInfo.Range("L10").Value = 1.182560277804791
$wsAltInfo.Range("L12").Value = 1.023100414119094
$wsAltInfo.Range("L13").Value = 1.30738877741298
$wsAltInfo.Range("L14").Value = 0.9389909155463747
$wsAltInfo.Range("L15").Value = 1.184877940454469
$wsAltInfo.Range("L18").Value = 1.282366543874901
$wsAltInfo.Range("L19").Value = 1.547292573912818
$wsAltInfo.Range("L20").Value = 1.152900358308349
$wsAltInfo.Range("L21").Value = 1.378578159221931
$wsAltInfo.Range("L23").Value = 1.28546803852881
$wsAltInfo.Range("L24").Value = 1.55

# alternatives_norm sheet - normalized C1.4 values (col D) recomputed
$wsAltNorm.Range("D2").Value = 0.00404375425268384
$wsAltNorm.Range("D3").Value = 0.004061857152499912
$wsAltNorm.Range("D4").Value = 0.005387262442215612
$wsAltNorm.Range("D5").Value = 0.000004218378453097304
$wsAltNorm.Range("D6").Value = 0.00626378845266073
$wsAltNorm.Range("D7").Value = 0.000004751231670574434
$wsAltNorm.Range("D8").Value = 0.000003572448176986866
$wsAltNorm.Range("D9").Value = 0.006739802024079785
$wsAltNorm.Range("D10").Value = 0.000003419491021794115
$wsAltNorm.Range("D11").Value = 0.008478847808840237
$wsAltNorm.Range("D12").Value = 0.000003952450997848124
$wsAltNorm.Range("D13").Value = 0.000003093000584482217
$wsAltNorm.Range("D14").Value = 0.000004306489217023876
$wsAltNorm.Range("D15").Value = 0.000003412802377882761
$wsAltNorm.Range("D16").Value = 0.0748843380126637
$wsAltNorm.Range("D17").Value = 0.4357372509970074
$wsAltNorm.Range("D18").Value = 0.000003153352894301896
$wsAltNorm.Range("D19").Value = 0.0000026134386740175
$wsAltNorm.Range("D20").Value = 0.000003507462048686707
$wsAltNorm.Range("D21").Value = 0.000002933278919032153
$wsAltNorm.Range("D22").Value = 0.4543544104139146
$wsAltNorm.Range("D23").Value = 0.000003145744687134988
$wsAltNorm.Range("D24").Value = 0.000002608873711408929

# result sheet - final AHP evaluation scores (col B) recomputed
$wsResult.Range("B2").Value = 0.1240858904054707
$wsResult.Range("B3").Value = 0.05909652950814774
$wsResult.Range("B4").Value = 0.05305799205377021
$wsResult.Range("B5").Value = 0.0488533249276431
$wsResult.Range("B6").Value = 0.04833357805990615
$wsResult.Range("B7").Value = 0.04747768887297976
$wsResult.Range("B8").Value = 0.04645410447782543
$wsResult.Range("B9").Value = 0.04604109166968942
$wsResult.Range("B10").Value = 0.0452419550870594
$wsResult.Range("B11").Value = 0.04033774628891104
$wsResult.Range("B12").Value = 0.03995918612765219
$wsResult.Range("B13").Value = 0.03915624901466493
$wsResult.Range("B14").Value = 0.03902574578025748
$wsResult.Range("B15").Value = 0.03789162019367157
$wsResult.Range("B16").Value = 0.03711900759619224
$wsResult.Range("B17").Value = 0.03698702851231435
$wsResult.Range("B18").Value = 0.03680355421438346
$wsResult.Range("B19").Value = 0.03542153714400613
$wsResult.Range("B20").Value = 0.02989540771484856
$wsResult.Range("B21").Value = 0.02904118697026819
$wsResult.Range("B22").Value = 0.02806178043104202
$wsResult.Range("B23").Value = 0.02646581161600669
$wsResult.Range("B24").Value = 0.02519198333328929
